# Updates cryptos list prices and 1h volume percentages (GitHub Actions data refresh).
# Column D ("Price") values are written with a leading apostrophe so Excel keeps
# them as literal text (matching the source XML's inlineStr cells) instead of
# auto-coercing dotted numeric-looking strings (e.g. "312.69") into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.523.49'
$ws.Range('E2').Value = '  +2.29%  '
$ws.Range('D3').Value = '''1.870.33'
$ws.Range('E3').Value = '  +1.39%  '
$ws.Range('E4').Value = '  +0.66%  '
$ws.Range('D5').Value = '''312.69'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('D7').Value = '''0.4778'
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('D8').Value = '''0.3774'
$ws.Range('E8').Value = '  +2.92%  '
$ws.Range('D9').Value = '''0.07357'
$ws.Range('E9').Value = '  +2.15%  '
$ws.Range('E10').Value = '  +1.15%  '
$ws.Range('D11').Value = '''20.72'
$ws.Range('E11').Value = '  +5.16%  '
$ws.Range('D12').Value = '''0.07851'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').Value = '''1.895.59'
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('D14').Value = '''5.444'
$ws.Range('E14').Value = '  +2.41%  '
$ws.Range('D15').Value = '''6.590'
$ws.Range('E15').Value = '  +2.86%  '
$ws.Range('D16').Value = '''90.78'
$ws.Range('E16').Value = '  +2.24%  '
$ws.Range('D17').Value = '''1.015'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('D18').Value = '''0.000008910'
$ws.Range('E18').Value = '  +3.12%  '
$ws.Range('D19').Value = '''1.012'
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D20').Value = '''14.96'
$ws.Range('E20').Value = '  +2.71%  '
$ws.Range('D21').Value = '''27.527.03'
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('D22').Value = '''5.140'
$ws.Range('E22').Value = '  +1.71%  '
$ws.Range('E23').Value = '  +1.01%  '
$ws.Range('D24').Value = '''1.962'
$ws.Range('E24').Value = '  +2.22%  '
$ws.Range('D25').Value = '''154.19'
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('E26').Value = '  +2.15%  '
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').Value = '''115.97'
$ws.Range('E28').Value = '  +1.50%  '
$ws.Range('D29').Value = '''4.999'
$ws.Range('E29').Value = '  +1.44%  '
$ws.Range('D30').Value = '''0.08930'
$ws.Range('E30').Value = '  +0.61%  '
$ws.Range('D31').Value = '''3.339'
$ws.Range('E31').Value = '  +0.92%  '
$ws.Range('E32').Value = '  +3.68%  '
$ws.Range('D33').Value = '''0.7559'
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('D34').Value = '''4.611'
$ws.Range('E34').Value = '  +2.81%  '
$ws.Range('D35').Value = '''2.689'
$ws.Range('E35').Value = '  -1.65%  '
$ws.Range('D36').Value = '''0.02052'
$ws.Range('E36').Value = '  +4.94%  '
$ws.Range('D37').Value = '''1.119'
$ws.Range('E37').Value = '  +2.31%  '
$ws.Range('D38').Value = '''0.05286'
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').Value = '''3.000'
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('D40').Value = '''0.5353'
$ws.Range('E40').Value = '  +2.90%  '
$ws.Range('D41').Value = '''7.077'
$ws.Range('E41').Value = '  +1.66%  '
$ws.Range('D42').Value = '''0.1526'
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('D43').Value = '''8.452'
$ws.Range('E43').Value = '  +2.77%  '
$ws.Range('D44').Value = '''10.66'
$ws.Range('E44').Value = '  +1.34%  '
$ws.Range('E45').Value = '  +1.97%  '
$ws.Range('D46').Value = '''1.013'
$ws.Range('E46').Value = '  +0.59%  '
$ws.Range('D47').Value = '''1.661'
$ws.Range('E47').Value = '  +3.66%  '
$ws.Range('D48').Value = '''102.78'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('D49').Value = '''67.40'
$ws.Range('E49').Value = '  +1.99%  '
$ws.Range('D50').Value = '''0.06090'
$ws.Range('E50').Value = '  +1.10%  '
$ws.Range('D51').Value = '''0.9293'
$ws.Range('E51').Value = '  +4.82%  '
